$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above the existing data rows (old row 2 "公共施設..."
# and old row 3 "水域情報..." shift down to become rows 5 and 6).
$ws.Rows("2:4").Insert()

# New sample rows (2-4)
$ws.Range("A2").Value = "サンプル"
$ws.Range("B2").Value = "NULL"
$ws.Range("C2").Value = "NULL"
$ws.Range("D2").Value = "#1f1f1f"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "〇"

$ws.Range("A3").Value = "施設サンプル"
$ws.Range("B3").Value = "NULL"
$ws.Range("C3").Value = "/public_facility"
$ws.Range("D3").Value = "#1f0000"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "〇"

$ws.Range("A4").Value = "地域施設サンプル"
$ws.Range("B4").Value = "112399_sakado_city"
$ws.Range("C4").Value = "NULL"
$ws.Range("D4").Value = "#1f1f00"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = "〇"

# Update the display-order values on the original rows, now at 5 and 6.
$ws.Range("E5").Value = 10
$ws.Range("E6").Value = 11

# Match the saved selection from the authored workbook.
$ws.Range("D8").Select() | Out-Null
